# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" column (E16:E29) listed the account-statement periods
# in ascending order (2310 .. 2411). The new periods (2401..2411) were
# added and the whole list is now shown in descending order
# (2411 .. 2310), newest period first.
#
# The "Valor Mora" column (F) carries the same value (46400) for every
# period except the oldest one (2310), which historically carried 21654.
# Because the oldest period moved from the first data row (16) to the
# last one (29), that odd value moves with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New, descending order of periods for rows 16 (newest) .. 29 (oldest).
$periodos = @("2411", "2410", "2409", "2408", "2407", "2406", "2405", "2404", "2403", "2402", "2401", "2312", "2311", "2310")

$firstRow = 16
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = $firstRow + $i
    $ws.Range("E$row").Value = $periodos[$i]
}

# The non-standard "Valor Mora" (21654, belonging to period 2310) now sits
# on the last data row instead of the first one; the rest keep 46400.
$ws.Range("F16").Value = 21654
$ws.Range("F29").Value = 46400
